# "upgrade left table until javakheti"
# For this municipality file (Oni):
#  1. Rename the (generically named "1") worksheet to the municipality name "Oni".
#  2. The table had a stray blank row (row 8) between the data rows (5-7) and the
#     footnote row (9). Delete that blank row so the footnote shifts up to row 8
#     and the used range becomes A1:AB8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Oni"
$ws.Rows(8).Delete()
